$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "Totaux" summary block (B15:F19, merged C:D) ---
$ws.Range("C15:D19").UnMerge()
$ws.Range("B15:F19").Clear()

# --- Rebuild the summary block next to the data table, at H3:L7 ---
$ws.Range("H3").Value = "Totaux"

$ws.Range("I3:J3").Merge()
$ws.Range("I3").Value = "Quantité vendue"
$ws.Range("I3").HorizontalAlignment = -4108

$ws.Range("K3").Value = "Montants HT"
$ws.Range("L3").Value = "Montants TTC"

# Dynamic-array UNIQUE formula spilling H4:H7
$ws.Range("H4:H7").FormulaArray = "=UNIQUE(B2:B11)"

$ws.Range("I4:J4").Merge()
$ws.Range("I4").HorizontalAlignment = -4108

$ws.Range("I5:J5").Merge()
$ws.Range("I5").HorizontalAlignment = -4108

$ws.Range("I6:J6").Merge()
$ws.Range("I6").HorizontalAlignment = -4108

$ws.Range("I7:J7").Merge()
$ws.Range("I7").HorizontalAlignment = -4108

# Column widths for the new K/L columns
$ws.Columns.Item(11).ColumnWidth = 13.7109375
$ws.Columns.Item(12).ColumnWidth = 13.28515625

# Keep the original "20%" value/format on I1 (percentage number format)
$ws.Range("I1").NumberFormat = "0%"

# Update selection to match the authored file
$ws.Range("H3:L7").Select()
